$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "52081"
$ws.Range("B2").Style = "Normal"
$ws.Range("E2").Value = "Уже есть"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "73297"
$ws.Range("B3").Style = "Normal"
$ws.Range("E3").Value = "Уже есть"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "40338"
$ws.Range("B4").Style = "Normal"
$ws.Range("E4").Value = "Уже есть"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "77209"
$ws.Range("B5").Style = "Normal"
$ws.Range("E5").Value = "Уже есть"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "38446"
$ws.Range("B6").Style = "Normal"
$ws.Range("E6").Value = "Уже есть"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "26754"
$ws.Range("B7").Style = "Normal"
$ws.Range("E7").Value = "Уже есть"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "57395"
$ws.Range("B8").Style = "Normal"
$ws.Range("E8").Value = "Уже есть"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "90364"
$ws.Range("B9").Style = "Normal"
$ws.Range("E9").Value = "Уже есть"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "34967"
$ws.Range("B10").Style = "Normal"
$ws.Range("E10").Value = "Уже есть"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "43969"
$ws.Range("B11").Style = "Normal"
$ws.Range("E11").Value = "Уже есть"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "40539"
$ws.Range("B12").Style = "Normal"
$ws.Range("E12").Value = "Уже есть"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "80525"
$ws.Range("B13").Style = "Normal"
$ws.Range("E13").Value = "Уже есть"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "64713"
$ws.Range("B14").Style = "Normal"
$ws.Range("E14").Value = "Уже есть"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "71808"
$ws.Range("B15").Style = "Normal"
$ws.Range("E15").Value = "Уже есть"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "51341"
$ws.Range("B16").Style = "Normal"
$ws.Range("E16").Value = "Уже есть"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "84158"
$ws.Range("B17").Style = "Normal"
$ws.Range("E17").Value = "Уже есть"

